# Disaggregation of commodity Copper
#
# 1) Rename the shared string "Copper ores and concentrates" -> "Copper".
#    Every worksheet carries this label in cell C4, so update it everywhere
#    the label occurs.
# 2) A handful of per-year sheets have their D4 total (the "Copper ores and
#    concentrates" / now "Copper" total) refreshed with a value that differs
#    from the previous one only in the last binary digit (re-aggregation
#    side effect of splitting out the new commodity).

$wb = $excel.ActiveWorkbook

$oldLabel = "Copper ores and concentrates"
$newLabel = "Copper"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("C4")
    if ($cell.Value() -eq $oldLabel) {
        $cell.Value = $newLabel
    }
}

# Sheet index -> refreshed D4 value (exact doubles from the target workbook)
$d4Updates = @{
    24 = 55579.97923991122
    26 = 64307.96100302236
    29 = 93395.27619719859
    33 = 136136.6974506026
    42 = 505872.9439998683
    49 = 1950034.592079028
    73 = 1953747.730931857
    86 = 1681427.682419382
    88 = 1854362.458575674
}

foreach ($sheetIndex in $d4Updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $ws.Range("D4").Value = $d4Updates[$sheetIndex]
}
